$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 714550.6
Write-Host "Cell H2 value set"
Write-Host $ws.Cells.Item(2, 8).Value
